$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 278

# Header for the new column K, matching the style of the existing header row (bold + centered)
$ws.Range("K1").Value = "PRODUCTO"
$ws.Range("K1").Font.Bold = $true
$ws.Range("K1").HorizontalAlignment = -4108

# Fill K2:K278 with "SOJA"
$ws.Range("K2:K$lastRow").Value = "SOJA"
